$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("full_signals - with decay")
$ws.Range("B16").Value = "unet_model_decay_coseno_with_normalization_zcore_6.keras"
Write-Output "done"
